$d = $word.ActiveDocument

# Locate the anchor paragraphs by their distinctive text content so the
# script does not depend on brittle absolute paragraph indices.
#
# $startIdx : paragraph "Варианты доработки программы:" - its TEXT is kept,
#             but its own paragraph mark (pilcrow) is removed, so it merges
#             with the paragraph that follows the deleted block.
# $lastFullIdx : paragraph "- запретить пользователю более одного раза
#             отправлять ответы на вопросы теста" - the last paragraph whose
#             text must be fully removed. The (now) empty paragraph right
#             after it is left in place on purpose: its paragraph mark is
#             the one that survives the merge performed in step 2 below.
$startIdx = -1
$lastFullIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -like "*Варианты доработки программы:*") {
        $startIdx = $i
    }
    if ($t -like "*запретить пользователю более одного раза отправлять ответы на вопросы теста*") {
        $lastFullIdx = $i
    }
}

if ($startIdx -eq -1 -or $lastFullIdx -eq -1) {
    throw "Could not locate the anchor paragraphs for the edit."
}

# Step 1: remove, in full (text + paragraph mark), every paragraph from the
# one right after "Варианты доработки программы:" through
# "- запретить пользователю ...теста" itself. The (now) empty paragraph
# that trails it is deliberately left alone for now: its paragraph mark is
# the one that needs to survive (see step 2).
$rangeStart = $d.Paragraphs.Item($startIdx + 1).Range.Start
$rangeEnd = $d.Paragraphs.Item($lastFullIdx).Range.End
$d.Range($rangeStart, $rangeEnd).Delete()

# Step 2: remove the paragraph mark of "Варианты доработки программы:"
# itself, merging it with the paragraph that now immediately follows
# (the formerly-empty paragraph). The merge is invisible since that
# paragraph has no text, and the combined paragraph simply precedes
# "Для пользователей можно:" afterwards.
$pStart = $d.Paragraphs.Item($startIdx)
$d.Range($pStart.Range.End - 1, $pStart.Range.End).Delete()
